$wb = $excel.ActiveWorkbook

# --- Sheet: Resource Utilization ---
$wsRU = $wb.Worksheets.Item("Resource Utilization")
$wsRU.Range("B2").Value = 3.77
$wsRU.Range("B3").Value = 1.22

# --- Sheet: Activity Times ---
$wsAT = $wb.Worksheets.Item("Activity Times")

# Row 2 - process summary row
$wsAT.Range("D2").Value = 4
$wsAT.Range("E2").Value = 15
$wsAT.Range("F2").Value = 164
$wsAT.Range("G2").Value = 87.5

# Row 3 - "Review AM using Asset Change Tracker (5.5.13.1)"
$wsAT.Range("C3").Value = 5
$wsAT.Range("D3").Value = 5
$wsAT.Range("E3").Value = 6
$wsAT.Range("F3").Value = 11
$wsAT.Range("G3").Value = 9.4

# Row 4 - now becomes "Work with REO RPO to Correct (5.5.13.3)"
$wsAT.Range("A4").Value = "Work with REO RPO to Correct (5.5.13.3)"
$wsAT.Range("B4").Value = "Activity Step"
$wsAT.Range("C4").Value = 2
$wsAT.Range("D4").Value = 2
$wsAT.Range("E4").Value = 129
$wsAT.Range("F4").Value = 150
$wsAT.Range("G4").Value = 139.5

# Row 5 - now becomes "Note Accuracy in Asset Change Tracker (5.5.13.2)"
$wsAT.Range("A5").Value = "Note Accuracy in Asset Change Tracker (5.5.13.2)"
$wsAT.Range("B5").Value = "Activity Step"
$wsAT.Range("C5").Value = 4
$wsAT.Range("D5").Value = 4
$wsAT.Range("E5").Value = 2
$wsAT.Range("F5").Value = 6
$wsAT.Range("G5").Value = 3.5

# Row 6 - now becomes "Create/Post Journal Entries (5.5.13.4)"
$wsAT.Range("A6").Value = "Create/Post Journal Entries (5.5.13.4)"
$wsAT.Range("B6").Value = "Stop"
$wsAT.Range("C6").Value = 4
$wsAT.Range("D6").Value = 4
$wsAT.Range("E6").Value = 4
$wsAT.Range("F6").Value = 6
$wsAT.Range("G6").Value = 5.25
